$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.269.51"
$ws.Range("E2").Value = "  +1.24%  "
$ws.Range("D3").Value = "2.992.86"
$ws.Range("E3").Value = "  +0.56%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'561.92"
$ws.Range("D6").Value = "'137.64"
$ws.Range("E6").Value = "  +4.19%  "
$ws.Range("E7").Value = "  -0.27%  "
$ws.Range("D8").Value = "'0.522"
$ws.Range("E8").Value = "  +0.56%  "
$ws.Range("D9").Value = "2.982.57"
$ws.Range("E9").Value = "  +0.62%  "
$ws.Range("E10").Value = "  +2.22%  "
$ws.Range("D11").Value = "'5.12"
$ws.Range("E11").Value = "  +5.58%  "
$ws.Range("E12").Value = "  +1.13%  "
$ws.Range("E13").Value = "  +1.64%  "
$ws.Range("D14").Value = "'33.61"
$ws.Range("E14").Value = "  +1.77%  "
$ws.Range("E15").Value = "  +1.65%  "
$ws.Range("D16").Value = "3.486.93"
$ws.Range("E16").Value = "  +0.53%  "
$ws.Range("E17").Value = "  +5.37%  "
$ws.Range("D18").Value = "2.991.87"
$ws.Range("E18").Value = "  +0.60%  "
$ws.Range("D19").Value = "59.224.64"
$ws.Range("E19").Value = "  +1.05%  "
$ws.Range("D20").Value = "'429.03"
$ws.Range("E20").Value = "  +1.56%  "
$ws.Range("D21").Value = "'13.67"
$ws.Range("E21").Value = "  +3.48%  "
$ws.Range("E22").Value = "  +4.50%  "
$ws.Range("D23").Value = "'7.11"
$ws.Range("E23").Value = "  +1.31%  "
$ws.Range("E24").Value = "  +1.81%  "
$ws.Range("D25").Value = "'80.96"
$ws.Range("E25").Value = "  +1.50%  "
$ws.Range("E26").Value = "  -0.03%  "
$ws.Range("E27").Value = "  +0.27%  "
$ws.Range("E28").Value = "  +6.85%  "
$ws.Range("E29").Value = "  +1.02%  "
$ws.Range("E30").Value = "  +1.58%  "
$ws.Range("D31").Value = "'25.72"
$ws.Range("E31").Value = "  +1.88%  "
$ws.Range("D32").Value = "'6.11"
$ws.Range("E32").Value = "  -1.60%  "
$ws.Range("D33").Value = "'0.0987"
$ws.Range("E33").Value = "  -6.21%  "
$ws.Range("B34").Value = "Filecoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D34").Value = "'5.92"
$ws.Range("E34").Value = "  +4.19%  "
$ws.Range("B35").Value = "Mantle"
$ws.Range("C35").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D35").Value = "'0.989"
$ws.Range("E35").Value = "  +4.31%  "
$ws.Range("D36").Value = "0.0₃0763"
$ws.Range("E36").Value = "  +9.86%  "
$ws.Range("E37").Value = "  -2.93%  "
$ws.Range("E38").Value = "  +0.85%  "
$ws.Range("D39").Value = "'8.69"
$ws.Range("E39").Value = "  +2.77%  "
$ws.Range("E40").Value = "  +2.40%  "
$ws.Range("D41").Value = "'401.31"
$ws.Range("E41").Value = "  +5.32%  "
$ws.Range("D42").Value = "'0.0351"
$ws.Range("E42").Value = "  -0.62%  "
$ws.Range("D43").Value = "2.755.22"
$ws.Range("E43").Value = "  +3.36%  "
$ws.Range("E44").Value = "  -1.33%  "
$ws.Range("E45").Value = "  +3.89%  "
$ws.Range("E46").Value = "  -0.02%  "
$ws.Range("D47").Value = "'34.73"
$ws.Range("E47").Value = "  +18.64%  "
$ws.Range("E48").Value = "  +0.77%  "
$ws.Range("D49").Value = "'121.15"
$ws.Range("E49").Value = "  -0.47%  "
$ws.Range("D50").Value = "'2.00"
$ws.Range("E50").Value = "  -0.18%  "
$ws.Range("D51").Value = "'23.41"
$ws.Range("E51").Value = "  -0.90%  "
